# Rotates the species-observation data (columns A,B,E,F,G,H,Q,R) across
# rows 2-5 by one position: new row2 = old row5, new row3 = old row2,
# new row4 = old row3, new row5 = old row4. All other columns in these
# rows are already identical across the block and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as one "record".
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

# Capture current (pre-edit) values for rows 2-5, one row at a time.
$rows = @(2, 3, 4, 5)
$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# New row r gets the old data from row (r - 1), wrapping row 2 <- row 5.
$sourceFor = @{ 2 = 5; 3 = 2; 4 = 3; 5 = 4 }

foreach ($r in $rows) {
    $src = $sourceFor[$r]
    $rowData = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $rowData[$c]
    }
}
